$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.770.17"
$ws.Range("E2").Value = '  -4.54%  '
$ws.Range("D3").Value = "'3.500.78"
$ws.Range("E3").Value = '  -5.20%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = "'605.22"
$ws.Range("E5").Value = '  -7.02%  '
$ws.Range("D6").Value = "'150.56"
$ws.Range("E6").Value = '  -6.82%  '
$ws.Range("D7").Value = "'3.503.72"
$ws.Range("E7").Value = '  -4.99%  '
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = "'0.480"
$ws.Range("E9").Value = '  -4.74%  '
$ws.Range("D10").Value = "'0.138"
$ws.Range("E10").Value = '  -5.63%  '
$ws.Range("D11").Value = "'6.87"
$ws.Range("E11").Value = '  -4.57%  '
$ws.Range("D12").Value = "'0.424"
$ws.Range("E12").Value = '  -4.94%  '
$ws.Range("D13").Value = "'0.0000218"
$ws.Range("E13").Value = '  -6.36%  '
$ws.Range("D14").Value = "'4.107.09"
$ws.Range("E14").Value = '  -4.80%  '
$ws.Range("D15").Value = "'31.25"
$ws.Range("E15").Value = '  -4.81%  '
$ws.Range("D16").Value = "'3.498.74"
$ws.Range("E16").Value = '  -4.76%  '
$ws.Range("D17").Value = "'66.865.75"
$ws.Range("E17").Value = '  -4.29%  '
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = "'6.29"
$ws.Range("E19").Value = '  -3.64%  '
$ws.Range("D20").Value = "'15.15"
$ws.Range("E20").Value = '  -6.16%  '
$ws.Range("D21").Value = "'442.64"
$ws.Range("E21").Value = '  -6.07%  '
$ws.Range("D22").Value = "'8.88"
$ws.Range("E22").Value = '  -16.03%  '
$ws.Range("D23").Value = "'0.622"
$ws.Range("E23").Value = '  -4.51%  '
$ws.Range("D24").Value = "'76.91"
$ws.Range("E24").Value = '  -3.80%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("D26").Value = "'3.655.67"
$ws.Range("E26").Value = '  -4.74%  '
$ws.Range("D27").Value = "'0.0000120"
$ws.Range("E27").Value = '  -4.86%  '
$ws.Range("D28").Value = "'10.05"
$ws.Range("E28").Value = '  -8.57%  '
$ws.Range("D29").Value = "'8.07"
$ws.Range("E29").Value = '  -12.19%  '
$ws.Range("D30").Value = "'2.51"
$ws.Range("E30").Value = '  -5.91%  '
$ws.Range("D31").Value = "'1.56"
$ws.Range("E31").Value = '  -9.70%  '
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("D33").Value = "'0.158"
$ws.Range("E33").Value = '  -4.11%  '
$ws.Range("D34").Value = "'25.44"
$ws.Range("E34").Value = '  -5.53%  '
$ws.Range("D35").Value = "'6.10"
$ws.Range("E35").Value = '  -7.34%  '
$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").Value = "'3.504.80"
$ws.Range("E36").Value = '  -4.96%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = "'1.85"
$ws.Range("E37").Value = '  -8.72%  '
$ws.Range("D38").Value = "'7.92"
$ws.Range("E38").Value = '  -6.86%  '
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("E40").Value = '  +0.17%  '
$ws.Range("D41").Value = "'172.13"
$ws.Range("E41").Value = '  -4.34%  '
$ws.Range("D42").Value = "'2.13"
$ws.Range("E42").Value = '  -6.43%  '
$ws.Range("D43").Value = "'5.48"
$ws.Range("E43").Value = '  -7.31%  '
$ws.Range("D44").Value = "'0.0853"
$ws.Range("E44").Value = '  -5.90%  '
$ws.Range("D45").Value = "'0.883"
$ws.Range("E45").Value = '  -4.93%  '
$ws.Range("D46").Value = "'45.13"
$ws.Range("D47").Value = "'26.81"
$ws.Range("E47").Value = '  -9.23%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").Value = "'1.21"
$ws.Range("E48").Value = '  -3.09%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = "'7.50"
$ws.Range("E49").Value = '  -4.63%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").Value = "'2.46"
$ws.Range("E50").Value = '  -10.36%  '
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = '  -6.93%  '
